# Agrega los metodos para obtener las normas en Matriz, VectorFila y VectorColumna.
# Llena las filas 8 y 9 de la hoja "Metricas" (antes vacias salvo formulas) con
# los datos de las nuevas tareas: "Metodos Norma 1, 2, inf en matriz" y
# "Metodos Norma 1, 2, inf en vector".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 8: Metodos Norma 1, 2, inf en matriz
$ws.Range("A8").Value = "Metodos Norma 1, 2, inf en matriz"
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = 24
$ws.Range("D8").Value = 0.017361111111111112
$ws.Range("E8").Value = 0.7125
$ws.Range("F8").Value = 0.72569444444444453
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.00069444444444444447

# Fila 9: Metodos Norma 1, 2, inf en vector
$ws.Range("A9").Value = "Metodos Norma 1, 2, inf en vector"
$ws.Range("B9").Value = 20
$ws.Range("D9").Value = 0.0069444444444444441
$ws.Range("E9").Value = 0.7284722222222223
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

# Refleja la nueva celda activa (selección) dejada por el autor tras editar.
$ws.Range("I10").Select() | Out-Null

$wb.Save()
